$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of data, matching the style of the existing empty B11 cell
$ws.Range("A11").Value = "channel_network_lines"
$ws.Range("B11").Value = "data/new_area/water_bodies_singleparts.gpkg"

# Move the active selection to B11
$ws.Range("B11").Select()
